{"js": "// Diff adds, right before the closing sectPr:\n//   1) a new empty paragraph (same \"en-US\" language formatting as the\n//      existing trailing empty paragraph)\n//   2) a new paragraph containing the run \"Vertical lines haben prior 0!\"\n//      (also \"en-US\")\n//\n// Word represents a brand-new, still-empty paragraph with a placeholder\n// run that carries the inherited character formatting forward (you can\n// see this any time you press Enter in Word and then look at the XML\n// before typing anything). We replicate that lifecycle here: insert both\n// paragraphs worth of text in one shot (using a newline to create the\n// paragraph break) and then drop the placeholder run's content from the\n// still-empty paragraph so it ends up with nothing but the paragraph\n// mark, matching a \"never typed into\" empty paragraph exactly.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// One insert creates both new paragraphs: the blank line, then the line\n// with our sentence (the embedded \"\\n\" becomes the paragraph break).\nlastParagraph.insertText(\"\\n\\nVertical lines haben prior 0!\", \"End\");\nawait context.sync();\n\n// Re-fetch the paragraph list now that two new paragraphs exist; the\n// newly inserted blank paragraph is the second-to-last one.\nconst updatedParagraphs = body.paragraphs;\nupdatedParagraphs.load(\"items\");\nawait context.sync();\n\nconst newBlankParagraph = updatedParagraphs.items[updatedParagraphs.items.length - 2];\nconst blankParagraphContent = newBlankParagraph.getRange(\"Content\");\nblankParagraphContent.delete();\nawait context.sync();\n", "ps1": "# Diff adds, right before the closing sectPr:\n#   1) a new empty paragraph (same \"en-US\" language formatting as the\n#      existing trailing empty paragraph)\n#   2) a new paragraph containing the run \"Vertical lines haben prior 0!\"\n#      (also \"en-US\")\n#\n# Word materializes a brand-new, still-empty paragraph with a placeholder\n# run that carries the inherited character formatting forward (the same\n# thing you'd see if you pressed Enter twice in Word and inspected the\n# XML before typing anything new). We reproduce that lifecycle: insert\n# both paragraph breaks plus the sentence in one go (two carriage returns\n# = two new paragraph marks), then delete the placeholder run's range\n# from the still-empty paragraph so only the paragraph mark is left,\n# matching an empty paragraph that was never typed into.\n\n$d = $word.ActiveDocument\n\n$lastParagraph = $d.Paragraphs($d.Paragraphs.Count)\n$endOfDoc = $lastParagraph.Range\n$endOfDoc.Collapse(0)  # wdCollapseEnd\n$endOfDoc.InsertAfter(\"`r`rVertical lines haben prior 0!\")\n\n$paragraphCount = $d.Paragraphs.Count\n$newBlankParagraph = $d.Paragraphs($paragraphCount - 1)\n$newBlankParagraph.Range.Delete()\n"}
